# 9.c.1 workbook update: add 2023 column (O) and refresh footnote source (МЦР КР).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the three footnote cells (row 8) with the new source text ---
$ws.Range("B8").Value = "*по данным МЦР КР"
$ws.Range("C8").Value = "*according to the MDD KR"
$ws.Range("A8").Value = "*КР СӨМ маалыматтары  боюнча"

# --- 2. Widen columns A:C from 35.7109375 to 38 characters ---
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 37.16666666666667

# --- 3. Add the new 2023 column (O) with its header + data, copying formats
#        from the adjacent 2022 column (N) so borders/number formats match ---
$ws.Range("N3:N7").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# --- 4. Reset the active selection to A1 (matches post-edit sheet state) ---
$ws.Range("A1").Select()
